$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($addr, $value) {
    # Force the cell to remain a text value even when it looks numeric
    # (e.g. "1.00", "0.380"), then restore the default "Normal" style so
    # no stray formatting is left behind on the cell.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-Text "D2" "58.398.28"
Set-Text "E2" "  +3.52%  "

# Row 3 - Ethereum
Set-Text "D3" "2.376.53"
Set-Text "E3" "  +1.92%  "

# Row 4 - TetherUSD
Set-Text "D4" "1.00"
Set-Text "E4" "  -0.02%  "

# Row 5 - BNB
Set-Text "D5" "545.19"
Set-Text "E5" "  +6.26%  "

# Row 6 - Solana
Set-Text "D6" "135.74"
Set-Text "E6" "  +2.24%  "

# Row 7 - USDC
Set-Text "D7" "1.00"
Set-Text "E7" "  -0.11%  "

# Row 8 - XRP
Set-Text "D8" "0.538"
Set-Text "E8" "  +1.03%  "

# Row 9 - LidoStakedEther
Set-Text "D9" "2.375.31"
Set-Text "E9" "  +1.76%  "

# Row 10 - Dogecoin
Set-Text "E10" "  +2.40%  "

# Row 11 - Toncoin
Set-Text "D11" "5.44"
Set-Text "E11" "  +2.83%  "

# Row 12 - TRON
Set-Text "D12" "0.154"
Set-Text "E12" "  +1.00%  "

# Row 13 - Cardano
Set-Text "D13" "0.356"
Set-Text "E13" "  +5.36%  "

# Row 14 - now Avalanche (was WrappedliquidstakedEther2.0)
Set-Text "B14" "Avalanche"
Set-Text "C14" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-Text "D14" "23.77"
Set-Text "E14" "  +0.64%  "

# Row 15 - now WrappedliquidstakedEther2.0 (was Avalanche)
Set-Text "B15" "WrappedliquidstakedEther2.0"
Set-Text "C15" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-Text "D15" "2.763.44"
Set-Text "E15" "  +0.57%  "

# Row 16 - WrappedBTC
Set-Text "D16" "58.247.64"
Set-Text "E16" "  +3.28%  "

# Row 17 - ShibaInu
Set-Text "E17" "  +1.35%  "

# Row 18 - WrappedEther
Set-Text "D18" "2.361.46"
Set-Text "E18" "  +0.69%  "

# Row 19 - BitcoinCash
Set-Text "D19" "339.50"
Set-Text "E19" "  +4.86%  "

# Row 20 - Chainlink
Set-Text "D20" "10.55"
Set-Text "E20" "  +1.52%  "

# Row 21 - Polkadot
Set-Text "D21" "4.25"
Set-Text "E21" "  +2.32%  "

# Row 22 - Uniswap
Set-Text "D22" "6.91"
Set-Text "E22" "  +4.22%  "

# Row 23 - Dai
Set-Text "E23" "  -0.16%  "

# Row 24 - Litecoin
Set-Text "D24" "62.44"
Set-Text "E24" "  +1.74%  "

# Row 25 - Kaspa
Set-Text "D25" "0.170"
Set-Text "E25" "  +4.05%  "

# Row 26 - InternetComputer(DFINITY)
Set-Text "D26" "8.56"
Set-Text "E26" "  -0.59%  "

# Row 27 - Binance-PegBSC-USD
Set-Text "D27" "0.993"
Set-Text "E27" "  -0.59%  "

# Row 28 - Fetch.AI
Set-Text "D28" "1.40"
Set-Text "E28" "  +8.03%  "

# Row 29 - now PancakeSwap (was Monero)
Set-Text "B29" "PancakeSwap"
Set-Text "C29" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-Text "D29" "1.77"
Set-Text "E29" "  +5.90%  "

# Row 30 - now Monero (was PancakeSwap)
Set-Text "B30" "Monero"
Set-Text "C30" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-Text "D30" "173.58"
Set-Text "E30" "  +3.65%  "

# Row 31 - PEPE
Set-Text "D31" "0.0₃0743"
Set-Text "E31" "  +3.09%  "

# Row 32 - Aptos
Set-Text "D32" "6.19"
Set-Text "E32" "  +1.37%  "

# Row 33 - EthereumClassic
Set-Text "D33" "18.61"
Set-Text "E33" "  +1.64%  "

# Row 34 - SuiNetwork
Set-Text "D34" "1.03"
Set-Text "E34" "  +16.39%  "

# Row 36 - FirstDigitalUSD
Set-Text "D36" "1.00"
Set-Text "E36" "  +0.20%  "

# Row 37 - ImmutableX
Set-Text "E37" "  +0.50%  "

# Row 38 - NEARProtocol
Set-Text "D38" "4.15"
Set-Text "E38" "  +5.12%  "

# Row 39 - Stacks
Set-Text "E39" "  +4.36%  "

# Row 40 - OKB
Set-Text "D40" "39.46"
Set-Text "E40" "  +2.73%  "

# Row 41 - Aave
Set-Text "D41" "150.55"
Set-Text "E41" "  +0.26%  "

# Row 42 - PolygonEcosystemToken
Set-Text "D42" "0.380"
Set-Text "E42" "  +1.50%  "

# Row 43 - Filecoin
Set-Text "D43" "3.65"
Set-Text "E43" "  +2.40%  "

# Row 44 - Bittensor
Set-Text "D44" "285.97"
Set-Text "E44" "  +3.02%  "

# Row 45 - InjectiveProtocol
Set-Text "D45" "19.34"
Set-Text "E45" "  +6.66%  "

# Row 46 - Stellar
Set-Text "D46" "0.0932"
Set-Text "E46" "  +0.76%  "

# Row 47 - Hedera
Set-Text "D47" "0.0507"
Set-Text "E47" "  +2.00%  "

# Row 48 - Mantle
Set-Text "D48" "0.563"
Set-Text "E48" "  +1.71%  "

# Row 49 - VeChain
Set-Text "E49" "  +2.38%  "

# Row 50 - EnergySwap
Set-Text "D50" "17.69"
Set-Text "E50" "  +4.04%  "

# Row 51 - Polygon
Set-Text "E51" "  +0.47%  "
